$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13, shifting rows 13:140 down to 14:141
$ws.Rows(13).Insert()

# Populate the newly inserted row 13 with the new record's data
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44545
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100112037
$ws.Range("G13").Value = "Cebollín"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 2960
$ws.Range("K13").Value = 900
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = 950
$ws.Range("N13").Value = "$/paquete 6 unidades"
$ws.Range("O13").Value = "Provincia del Elquí"
$ws.Range("P13").Value = 158
$ws.Range("Q13").Value = 6
$ws.Range("R13").Value = "Hortaliza"
